# Update "want to go" counts (column F) across sheets to match the
# latest scrape output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 663
$wsExhibit.Range("F8").Value = 2594
$wsExhibit.Range("F9").Value = 4138

# Sheet "演出" (Performance)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 60

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 663
$wsAll.Range("F8").Value = 2594
$wsAll.Range("F9").Value = 4138
$wsAll.Range("F11").Value = 60
